$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear values in column D (Attendees) for rows that still keep their cell/style,
# and fully clear (remove) the cells for rows whose <c> element disappears.
$ws.Range("D3:D6").ClearContents()
$ws.Range("D8:D11").Clear()
$ws.Range("D13:D16").Clear()
$ws.Range("D18:D24").Clear()
$ws.Range("D26:D28").Clear()
$ws.Range("D30:D31").Clear()
$ws.Range("D33:D34").Clear()

# Scroll the view back to the top and select D3 (matches the reset view state).
$ws.Range("D3").Select()
